$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds "K" (strikeouts), regenerated to replace the old "Strike#" values.
$kValues = @{
    2  = 6
    3  = 0
    4  = 1
    5  = 4
    6  = 1
    7  = 1
    8  = 2
    9  = 2
    10 = 1
    11 = 1
    12 = 4
    13 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
